# PTW-TimeSheet: fill in actual start/end times and task details for
# 28-04-2022 sheet, rows 92-103 (the timesheet entries that were still
# blank / placeholder), and leave the cursor near the last edited cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 92 - adjust the start/end time of the first entry of the day
$ws.Range("D92").Value = 0.36458333333333331
$ws.Range("E92").Value = 0.375

# Row 93
$ws.Range("B93").Value = "Worked on Web api(SOC,Factory pattern)"
$ws.Range("C93").Value = "Project"
$ws.Range("D93").Value = 0.375
$ws.Range("E93").Value = 0.45833333333333331

# Row 94
$ws.Range("B94").Value = "Morning Break"
$ws.Range("C94").Value = "Lunch and Break"
$ws.Range("D94").Value = 0.46527777777777773
$ws.Range("E94").Value = 0.47916666666666669

# Row 95
$ws.Range("B95").Value = "Worked on Web api"
$ws.Range("C95").Value = "Project"
$ws.Range("D95").Value = 0.47916666666666669
$ws.Range("E95").Value = 0.54166666666666663

# Row 96
$ws.Range("B96").Value = "Customer Meeting"
$ws.Range("C96").Value = "Project"
$ws.Range("D96").Value = 0.54166666666666663
$ws.Range("E96").Value = 0.57291666666666663

# Row 97
$ws.Range("B97").Value = "Lunch Break"
$ws.Range("C97").Value = "Lunch and Break"
$ws.Range("D97").Value = 0.57291666666666663
$ws.Range("E97").Value = 0.59375

# Row 98
$ws.Range("B98").Value = "Logging"
$ws.Range("C98").Value = "Exploration "
$ws.Range("D98").Value = 0.59722222222222221
$ws.Range("E98").Value = 0.625

# Row 99
$ws.Range("B99").Value = "Tried Console logging"
$ws.Range("C99").Value = "Non Project"
$ws.Range("D99").Value = 0.625
$ws.Range("E99").Value = 0.64583333333333337

# Row 100
$ws.Range("B100").Value = "TypeScript Session"
$ws.Range("C100").Value = "Exploration "
$ws.Range("D100").Value = 0.64583333333333337
$ws.Range("E100").Value = 0.70833333333333337

# Row 101 (previously blank placeholder row)
$ws.Range("B101").Value = "Evening Break"
$ws.Range("C101").Value = "Lunch and Break"
$ws.Range("D101").Value = 0.70833333333333337
$ws.Range("E101").Value = 0.72222222222222221

# Row 102 (previously blank placeholder row)
$ws.Range("B102").Value = "Logging"
$ws.Range("C102").Value = "Project"
$ws.Range("D102").Value = 0.72222222222222221
$ws.Range("E102").Value = 0.74305555555555547

# Row 103 (previously blank placeholder row)
$ws.Range("B103").Value = "TypeScript and Entity Framework"
$ws.Range("C103").Value = "Exploration "
$ws.Range("D103").Value = 0.80208333333333337
$ws.Range("E103").Value = 0.83333333333333337

# Reflect where the author ended up working (scrolled down, cursor on F103)
$ws.Range("F103").Select()
$excel.ActiveWindow.ScrollRow = 92
